$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 4 / Row 5 header relabeling + new header block (Z:AC mirrors old V:Y)
# ---------------------------------------------------------------------------
# V4 used to hold "Standard normalisation" (shared string 32); it now holds
# the new label "Standard normalisation for each sample", and the old label
# moves to Z4.
$ws.Range("V4").Value = "Standard normalisation for each sample"
$ws.Range("Z4").Value = "Standard normalisation"

# V5 used to hold "F14C " (with trailing space); it now holds the new,
# trimmed label "F14C", and the old header row (F14C / delta_F14C / Time T /
# delta_T) is duplicated into Z5:AC5.
$ws.Range("V5").Value = "F14C"
$ws.Range("W5").Value = "delta_F14C"
$ws.Range("X5").Value = "Time T"
$ws.Range("Y5").Value = "delta_T"

$ws.Range("Z5").Value = "F14C "
$ws.Range("AA5").Value = "delta_F14C"
$ws.Range("AB5").Value = "Time T"
$ws.Range("AC5").Value = "delta_T"

# ---------------------------------------------------------------------------
# Row 11 (blank 1 of 2): new per-sample normalisation in V:Y, the previous
# pooled-blank normalisation formulas move (unchanged) to Z:AC.
# ---------------------------------------------------------------------------
$ws.Range("Z11").Formula = "=`$B`$44*1.34066/`$B`$41"
$ws.Range("AA11").Formula = "=Z11*SQRT((U11/T11)^2+(U12/T12)^2+`$D`$42/(`$B`$41^2))"
$ws.Range("AB11").Formula = "=-8033*LN(Z11)"
$ws.Range("AC11").Formula = "=8033*AA11/Z11"

$ws.Range("V11").Formula = "=T11*1.34066/`$B`$41"
$ws.Range("W11").Formula = "=V11*SQRT((U11/T11)^2+`$D`$42/(`$B`$41^2))"
$ws.Range("X11").Formula = "=-8033*LN(V11)"
$ws.Range("Y11").Formula = "=8033*W11/V11"

# Row 12 (blank 2 of 2): only gets the new per-sample normalisation columns.
$ws.Range("V12").Formula = "=T12*1.34066/`$B`$41"
$ws.Range("W12").Formula = "=V12*SQRT((U12/T12)^2+`$D`$42/(`$B`$41^2))"
$ws.Range("X12").Formula = "=-8033*LN(V12)"
$ws.Range("Y12").Formula = "=8033*W12/V12"

# ---------------------------------------------------------------------------
# Row 19 (sample 1 of 15): new per-sample normalisation in V:Y, the previous
# pooled-sample normalisation formulas move (unchanged) to Z:AC.
# ---------------------------------------------------------------------------
$ws.Range("Z19").Formula = "=`$B`$45*1.34066/`$B`$41"
$ws.Range("AA19").Formula = "=Z19*SQRT((U19/T19)^2+(U20/T20)^2+(U21/T21)^2+(U22/T22)^2+(U23/T23)^2+(U24/T24)^2+(U25/T25)^2+(U26/T26)^2+(U27/T27)^2+(U28/T28)^2+(U29/T29)^2+(U30/T30)^2+(U31/T31)^2+(U32/T32)^2+(U33/T33)^2+`$D`$42/(`$B`$41^2))"
$ws.Range("AB19").Formula = "=-8033*LN(Z19)"
$ws.Range("AC19").Formula = "=8033*AA19/Z19"

$ws.Range("V19").Formula = "=T19*1.34066/`$B`$41"
$ws.Range("W19").Formula = "=V19*SQRT((U19/T19)^2+`$D`$42/(`$B`$41^2))"
$ws.Range("X19").Formula = "=-8033*LN(V19)"
$ws.Range("Y19").Formula = "=8033*W19/V19"

# Rows 20-33 (samples 2-15 of 15): only the new per-sample normalisation
# columns.
for ($r = 20; $r -le 33; $r++) {
    $ws.Range("V$r").Formula = "=T$r*1.34066/`$B`$41"
    $ws.Range("W$r").Formula = "=V$r*SQRT((U$r/T$r)^2+`$D`$42/(`$B`$41^2))"
    $ws.Range("X$r").Formula = "=-8033*LN(V$r)"
    $ws.Range("Y$r").Formula = "=8033*W$r/V$r"
}

# ---------------------------------------------------------------------------
# Sheet view: the active selection moved from Q11 to V19:V33 while scrolled
# so row 15 / column R sits at the top-left of the viewport.
# ---------------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 15
$win.ScrollColumn = 18
$ws.Range("V19:V33").Select()
